$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells - copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I and J
$values = @{
    2 = @(6, 8)
    3 = @(7, 8)
    4 = @(2, 6)
    5 = @(9, 9)
    6 = @(7, 8)
    7 = @(8, 9)
    8 = @(6, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
